# Upgrade the left table with the 2023 column (Dusheti municipality sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dusheti")

# New year header + the three data rows that follow the existing B:J pattern.
$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 1078.3
$ws.Range("K5").Value = 460.4
$ws.Range("K6").Value = 1454.5

# Carry the formatting (number format / borders / alignment) from the
# previous year's column (J) onto the newly added column (K).
$ws.Range("J3:J6").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122)
$excel.CutCopyMode = $false
